# Insert a new weekly price-report row for "Agrícola del Norte S.A. de Arica"
# (Manzana / Royal Gala) above the existing row 235, shifting the rest of the
# table down by one row (old row 235 -> new row 236, ..., old row 277 -> new
# row 278). The worksheet's used-range dimension grows from A1:T277 to
# A1:T278 automatically once the new row contains data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 235:277 down to 236:278, leaving a blank row 235 to fill in.
$ws.Rows.Item(235).Insert()

$ws.Range("A235").Value = 1
$ws.Range("B235").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C235").Value = "Arica y Parinacota"
$ws.Range("D235").Value = 45166
$ws.Range("E235").Value = 15
$ws.Range("F235").Value = "Fruta"
$ws.Range("G235").Value = 100104
$ws.Range("H235").Value = "Frutos de pepita"
$ws.Range("I235").Value = 100104002
$ws.Range("J235").Value = "Manzana"
$ws.Range("K235").Value = "Royal Gala"
$ws.Range("L235").Value = "Calibre 80"
$ws.Range("M235").Value = 300
$ws.Range("N235").Value = 23000
$ws.Range("O235").Value = 25000
$ws.Range("P235").Value = 24333
$ws.Range("Q235").Value = "$/caja 18 kilos embalada"
$ws.Range("R235").Value = "Provincia de Cachapoal"
$ws.Range("S235").Value = 1352
$ws.Range("T235").Value = 18
